$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Adjust column widths (closest achievable values through the
#     character-width/pixel quantized ColumnWidth API) for columns C:L ---
$ws.Columns.Item(3).ColumnWidth  = 14.333333333333334
$ws.Columns.Item(4).ColumnWidth  = 12.833333333333334
$ws.Columns.Item(5).ColumnWidth  = 18.5
$ws.Columns.Item(6).ColumnWidth  = 18.666666666666668
$ws.Columns.Item(7).ColumnWidth  = 28.333333333333332
$ws.Columns.Item(8).ColumnWidth  = 28.833333333333332
$ws.Columns.Item(9).ColumnWidth  = 18.0
$ws.Columns.Item(10).ColumnWidth = 21.833333333333332
$ws.Columns.Item(11).ColumnWidth = 24.333333333333332
$ws.Columns.Item(12).ColumnWidth = 20.333333333333332

# --- Apply the HH:MM:SS time format to the new data row's time columns ---
$ws.Range("E3:L3").NumberFormat = "HH:MM:SS"

# --- Add the new data row (row 3) ---
$ws.Range("A3").Value = 11
$ws.Range("B3").Value = "2b"
$ws.Range("C3").Value = 3145
$ws.Range("D3").Value = "lucas"
$ws.Range("E3").Value = 0.555555555555556
$ws.Range("F3").Value = 0.541666666666667
$ws.Range("G3").Value = 0.791666666666667
$ws.Range("H3").Value = 0.833333333333333
$ws.Range("I3").Value = 0.809027777777778
$ws.Range("J3").Value = 0.822916666666667
$ws.Range("K3").Value = 0.875
$ws.Range("L3").Value = 0.916666666666667

# --- Move the active selection to H6 (matches the saved selection state) ---
[void]$ws.Range("H6").Select()
